$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 19.54007855928535
$ws.Cells.Item(2, 3).Value = 6.403381025117031
$ws.Cells.Item(2, 4).Value = 4.736560110515995
$ws.Cells.Item(2, 5).Value = 10.65486524604365
$ws.Cells.Item(2, 6).Value = 50.32538137073504
$ws.Cells.Item(2, 8).Value = 7.344005520526261
$ws.Cells.Item(2, 9).Value = 36.22301452445404
$ws.Cells.Item(2, 10).Value = 9.964561230390371
$ws.Cells.Item(2, 11).Value = 17.35877020247761
$ws.Cells.Item(2, 12).Value = 12.05623818284223
$ws.Cells.Item(2, 14).Value = 24.47797805936214
$ws.Cells.Item(3, 2).Value = 19.42173598948603
$ws.Cells.Item(3, 3).Value = 6.274015419379142
$ws.Cells.Item(3, 4).Value = 4.73528851569362
$ws.Cells.Item(3, 5).Value = 10.67026975130614
$ws.Cells.Item(3, 6).Value = 50.34835583333911
$ws.Cells.Item(3, 8).Value = 7.344005520526261
$ws.Cells.Item(3, 9).Value = 36.27675806777497
$ws.Cells.Item(3, 10).Value = 9.981215834887967
$ws.Cells.Item(3, 11).Value = 17.28004288095998
$ws.Cells.Item(3, 12).Value = 12.06628163864889
$ws.Cells.Item(3, 14).Value = 24.53379105072665
$ws.Cells.Item(4, 2).Value = 19.35328567392257
$ws.Cells.Item(4, 3).Value = 6.195227297311513
$ws.Cells.Item(4, 4).Value = 4.735032054012545
$ws.Cells.Item(4, 5).Value = 10.68076458328744
$ws.Cells.Item(4, 6).Value = 50.37157890864938
$ws.Cells.Item(4, 8).Value = 7.344005520526261
$ws.Cells.Item(4, 9).Value = 36.315312857768
$ws.Cells.Item(4, 10).Value = 9.992088785473467
$ws.Cells.Item(4, 11).Value = 17.23540360998681
$ws.Cells.Item(4, 12).Value = 12.07415306756759
$ws.Cells.Item(4, 14).Value = 24.56999178885823
$ws.Cells.Item(5, 2).Value = 19.32647562647824
$ws.Cells.Item(5, 3).Value = 6.163331557302289
$ws.Cells.Item(5, 4).Value = 4.73506010751039
$ws.Cells.Item(5, 5).Value = 10.68530240393084
$ws.Cells.Item(5, 6).Value = 50.3833341516757
$ws.Cells.Item(5, 8).Value = 7.344005520526261
$ws.Cells.Item(5, 9).Value = 36.33242004420122
$ws.Cells.Item(5, 10).Value = 9.996682704421696
$ws.Cells.Item(5, 11).Value = 17.21815753813513
$ws.Cells.Item(5, 12).Value = 12.07778988433805
$ws.Cells.Item(5, 14).Value = 24.58523027390848
$ws.Cells.Item(6, 2).Value = 19.3220899991836
$ws.Cells.Item(6, 3).Value = 6.15804949284908
$ws.Cells.Item(6, 4).Value = 4.735072792810843
$ws.Cells.Item(6, 5).Value = 10.6860716887607
$ws.Cells.Item(6, 6).Value = 50.38542447127992
$ws.Cells.Item(6, 8).Value = 7.344005520526261
$ws.Cells.Item(6, 9).Value = 36.3353449398677
$ws.Cells.Item(6, 10).Value = 9.997455384599574
$ws.Cells.Item(6, 11).Value = 17.21535130294266
$ws.Cells.Item(6, 12).Value = 12.07841970682214
$ws.Cells.Item(6, 14).Value = 24.58779001148926
$ws.Cells.Item(7, 2).Value = 19.35291968374514
$ws.Cells.Item(7, 3).Value = 6.194796220061314
$ws.Cells.Item(7, 4).Value = 4.735031894662532
$ws.Cells.Item(7, 5).Value = 10.68082472413887
$ws.Cells.Item(7, 6).Value = 50.37172816692188
$ws.Cells.Item(7, 8).Value = 7.344005520526261
$ws.Cells.Item(7, 9).Value = 36.31553792168611
$ws.Cells.Item(7, 10).Value = 9.99215007971301
$ws.Cells.Item(7, 11).Value = 17.23516717981099
$ws.Cells.Item(7, 12).Value = 12.0742003767803
$ws.Cells.Item(7, 14).Value = 24.5701953300658
$ws.Cells.Item(8, 2).Value = 19.4984151000871
$ws.Cells.Item(8, 3).Value = 6.358671809577812
$ws.Cells.Item(8, 4).Value = 4.736013288292685
$ws.Cells.Item(8, 5).Value = 10.6599619036836
$ws.Cells.Item(8, 6).Value = 50.33141054328737
$ws.Cells.Item(8, 8).Value = 7.344005520526261
$ws.Cells.Item(8, 9).Value = 36.24039140489471
$ws.Cells.Item(8, 10).Value = 9.970169715456919
$ws.Cells.Item(8, 11).Value = 17.33086623704553
$ws.Cells.Item(8, 12).Value = 12.05934766689085
$ws.Cells.Item(8, 14).Value = 24.49682197577688
$ws.Cells.Item(9, 2).Value = 19.81598962260661
$ws.Cells.Item(9, 3).Value = 6.682997188148645
$ws.Cells.Item(9, 4).Value = 4.742065288817394
$ws.Cells.Item(9, 5).Value = 10.62725225823489
$ws.Cells.Item(9, 6).Value = 50.32469321505917
$ws.Cells.Item(9, 8).Value = 7.344005520526261
$ws.Cells.Item(9, 9).Value = 36.1371680047985
$ws.Cells.Item(9, 10).Value = 9.932180642495519
$ws.Cells.Item(9, 11).Value = 17.54719325208677
$ws.Cells.Item(9, 12).Value = 12.04372424308827
$ws.Cells.Item(9, 14).Value = 24.3682272436387
$ws.Cells.Item(10, 2).Value = 20.06731244532049
$ws.Cells.Item(10, 3).Value = 6.920312362171358
$ws.Cells.Item(10, 4).Value = 4.748982927110498
$ws.Cells.Item(10, 5).Value = 10.60819232959844
$ws.Cells.Item(10, 6).Value = 50.36383344612475
$ws.Cells.Item(10, 8).Value = 7.344005520526261
$ws.Cells.Item(10, 9).Value = 36.0883046222967
$ws.Cells.Item(10, 10).Value = 9.907361756555954
$ws.Cells.Item(10, 11).Value = 17.72259546934967
$ws.Cells.Item(10, 12).Value = 12.04044265439576
$ws.Cells.Item(10, 14).Value = 24.2830218786368
$ws.Cells.Item(11, 2).Value = 20.1851728933797
$ws.Cells.Item(11, 3).Value = 7.027509599922821
$ws.Cells.Item(11, 4).Value = 4.752656299845672
$ws.Cells.Item(11, 5).Value = 10.60059485629848
$ws.Cells.Item(11, 6).Value = 50.39118621767945
$ws.Cells.Item(11, 8).Value = 7.344005520526261
$ws.Cells.Item(11, 9).Value = 36.07194284718362
$ws.Cells.Item(11, 10).Value = 9.896736805067983
$ws.Cells.Item(11, 11).Value = 17.80572971091559
$ws.Cells.Item(11, 12).Value = 12.04072051087904
$ws.Cells.Item(11, 14).Value = 24.24626333430215
$ws.Cells.Item(12, 2).Value = 20.23027788076839
$ws.Cells.Item(12, 3).Value = 7.067950097207169
$ws.Cells.Item(12, 4).Value = 4.754122087003588
$ws.Cells.Item(12, 5).Value = 10.59787165453562
$ws.Cells.Item(12, 6).Value = 50.40291318713778
$ws.Cells.Item(12, 8).Value = 7.344005520526261
$ws.Cells.Item(12, 9).Value = 36.06659100947495
$ws.Cells.Item(12, 10).Value = 9.892808656703862
$ws.Cells.Item(12, 11).Value = 17.8376698372222
$ws.Cells.Item(12, 12).Value = 12.04107938619877
$ws.Cells.Item(12, 14).Value = 24.23263099853331
$ws.Cells.Item(13, 2).Value = 20.22054319497168
$ws.Cells.Item(13, 3).Value = 7.059247976453714
$ws.Cells.Item(13, 4).Value = 4.753803094011258
$ws.Cells.Item(13, 5).Value = 10.59845131265779
$ws.Cells.Item(13, 6).Value = 50.40032675650156
$ws.Cells.Item(13, 8).Value = 7.344005520526261
$ws.Cells.Item(13, 9).Value = 36.0677060817037
$ws.Cells.Item(13, 10).Value = 9.893650421644743
$ws.Cells.Item(13, 11).Value = 17.83077089980878
$ws.Cells.Item(13, 12).Value = 12.04099083089654
$ws.Cells.Item(13, 14).Value = 24.23555419521947
$ws.Cells.Item(14, 2).Value = 20.18887441082535
$ws.Cells.Item(14, 3).Value = 7.030839967208652
$ws.Cells.Item(14, 4).Value = 4.752775398286557
$ws.Cells.Item(14, 5).Value = 10.60036773742702
$ws.Cells.Item(14, 6).Value = 50.39212358990679
$ws.Cells.Item(14, 8).Value = 7.344005520526261
$ws.Cells.Item(14, 9).Value = 36.0714856314776
$ws.Cells.Item(14, 10).Value = 9.896411726327022
$ws.Cells.Item(14, 11).Value = 17.80834837260824
$ws.Cells.Item(14, 12).Value = 12.04074495796863
$ws.Cells.Item(14, 14).Value = 24.24513603833628
$ws.Cells.Item(15, 2).Value = 20.16953704035936
$ws.Cells.Item(15, 3).Value = 7.013418091400537
$ws.Cells.Item(15, 4).Value = 4.75215561223987
$ws.Cells.Item(15, 5).Value = 10.60156161683077
$ws.Cells.Item(15, 6).Value = 50.38727707546718
$ws.Cells.Item(15, 8).Value = 7.344005520526261
$ws.Cells.Item(15, 9).Value = 36.07391063843585
$ws.Cells.Item(15, 10).Value = 9.898115503644558
$ws.Cells.Item(15, 11).Value = 17.79467304371473
$ws.Cells.Item(15, 12).Value = 12.04062735707664
$ws.Cells.Item(15, 14).Value = 24.25104259700313
$ws.Cells.Item(16, 2).Value = 20.05967816775961
$ws.Cells.Item(16, 3).Value = 6.913287715876085
$ws.Cells.Item(16, 4).Value = 4.748753377754018
$ws.Cells.Item(16, 5).Value = 10.60871039126748
$ws.Cells.Item(16, 6).Value = 50.362237692119
$ws.Cells.Item(16, 8).Value = 7.344005520526261
$ws.Cells.Item(16, 9).Value = 36.08949199240327
$ws.Cells.Item(16, 10).Value = 9.908069477043295
$ws.Cells.Item(16, 11).Value = 17.7172278814688
$ws.Cells.Item(16, 12).Value = 12.04046004614721
$ws.Cells.Item(16, 14).Value = 24.28546437163719
$ws.Cells.Item(17, 2).Value = 19.99316400978027
$ws.Cells.Item(17, 3).Value = 6.851634830275172
$ws.Cells.Item(17, 4).Value = 4.746800402513359
$ws.Cells.Item(17, 5).Value = 10.61337039939147
$ws.Cells.Item(17, 6).Value = 50.34931977350787
$ws.Cells.Item(17, 8).Value = 7.344005520526261
$ws.Cells.Item(17, 9).Value = 36.10055356297338
$ws.Cells.Item(17, 10).Value = 9.914346043261665
$ws.Cells.Item(17, 11).Value = 17.67055879427296
$ws.Cells.Item(17, 12).Value = 12.04081026324036
$ws.Cells.Item(17, 14).Value = 24.30709334239482
$ws.Cells.Item(18, 2).Value = 19.95524141836927
$ws.Cells.Item(18, 3).Value = 6.816104780693835
$ws.Cells.Item(18, 4).Value = 4.745726718975011
$ws.Cells.Item(18, 5).Value = 10.61615173757421
$ws.Cells.Item(18, 6).Value = 50.34278869305906
$ws.Cells.Item(18, 8).Value = 7.344005520526261
$ws.Cells.Item(18, 9).Value = 36.10746805119955
$ws.Cells.Item(18, 10).Value = 9.9180188001391
$ws.Cells.Item(18, 11).Value = 17.64403197321807
$ws.Cells.Item(18, 12).Value = 12.04117842632718
$ws.Cells.Item(18, 14).Value = 24.31972220385935
$ws.Cells.Item(19, 2).Value = 19.94245997449477
$ws.Cells.Item(19, 3).Value = 6.804064446268725
$ws.Cells.Item(19, 4).Value = 4.745371738171398
$ws.Cells.Item(19, 5).Value = 10.61711081695998
$ws.Cells.Item(19, 6).Value = 50.34073188226614
$ws.Cells.Item(19, 8).Value = 7.344005520526261
$ws.Cells.Item(19, 9).Value = 36.1099039960228
$ws.Cells.Item(19, 10).Value = 9.919273102879529
$ws.Cells.Item(19, 11).Value = 17.63510536935941
$ws.Cells.Item(19, 12).Value = 12.04133174334085
$ws.Cells.Item(19, 14).Value = 24.32403050224852
$ws.Cells.Item(20, 2).Value = 20.00021017119133
$ws.Cells.Item(20, 3).Value = 6.858205344400675
$ws.Cells.Item(20, 4).Value = 4.747003171757724
$ws.Cells.Item(20, 5).Value = 10.61286388146067
$ws.Cells.Item(20, 6).Value = 50.35060189353521
$ws.Cells.Item(20, 8).Value = 7.344005520526261
$ws.Cells.Item(20, 9).Value = 36.09931888890932
$ws.Cells.Item(20, 10).Value = 9.913671411600927
$ws.Cells.Item(20, 11).Value = 17.67549424776278
$ws.Cells.Item(20, 12).Value = 12.04075573302329
$ws.Cells.Item(20, 14).Value = 24.30477140157507
$ws.Cells.Item(21, 2).Value = 20.19816371798312
$ws.Cells.Item(21, 3).Value = 7.039188570256455
$ws.Cells.Item(21, 4).Value = 4.753075236193554
$ws.Cells.Item(21, 5).Value = 10.59980066704343
$ws.Cells.Item(21, 6).Value = 50.39449593822678
$ws.Cells.Item(21, 8).Value = 7.344005520526261
$ws.Cells.Item(21, 9).Value = 36.07035257770094
$ws.Cells.Item(21, 10).Value = 9.895598081045016
$ws.Cells.Item(21, 11).Value = 17.81492213240166
$ws.Cells.Item(21, 12).Value = 12.0408103007917
$ws.Cells.Item(21, 14).Value = 24.24231382489119
$ws.Cells.Item(22, 2).Value = 20.33028306839981
$ws.Cells.Item(22, 3).Value = 7.156562603624842
$ws.Cells.Item(22, 4).Value = 4.757478978573962
$ws.Cells.Item(22, 5).Value = 10.59215930631872
$ws.Cells.Item(22, 6).Value = 50.43116074817979
$ws.Cells.Item(22, 8).Value = 7.344005520526261
$ws.Cells.Item(22, 9).Value = 36.05634098336421
$ws.Cells.Item(22, 10).Value = 9.884341366441765
$ws.Cells.Item(22, 11).Value = 17.90871106935058
$ws.Cells.Item(22, 12).Value = 12.04232395488613
$ws.Cells.Item(22, 14).Value = 24.20316871881339
$ws.Cells.Item(23, 2).Value = 20.25952877586863
$ws.Cells.Item(23, 3).Value = 7.094014835495637
$ws.Cells.Item(23, 4).Value = 4.755089111370412
$ws.Cells.Item(23, 5).Value = 10.59615580867614
$ws.Cells.Item(23, 6).Value = 50.41086363071553
$ws.Cells.Item(23, 8).Value = 7.344005520526261
$ws.Cells.Item(23, 9).Value = 36.06336901854225
$ws.Cells.Item(23, 10).Value = 9.890298604907736
$ws.Cells.Item(23, 11).Value = 17.85841769647774
$ws.Cells.Item(23, 12).Value = 12.04138119667442
$ws.Cells.Item(23, 14).Value = 24.22390814899576
$ws.Cells.Item(24, 2).Value = 19.99702361068591
$ws.Cells.Item(24, 3).Value = 6.85523507632498
$ws.Cells.Item(24, 4).Value = 4.746911346765049
$ws.Cells.Item(24, 5).Value = 10.61309255973064
$ws.Cells.Item(24, 6).Value = 50.35001945691562
$ws.Cells.Item(24, 8).Value = 7.344005520526261
$ws.Cells.Item(24, 9).Value = 36.09987535606204
$ws.Cells.Item(24, 10).Value = 9.913976212330835
$ws.Cells.Item(24, 11).Value = 17.67326198107839
$ws.Cells.Item(24, 12).Value = 12.04077986643759
$ws.Cells.Item(24, 14).Value = 24.30582054637216
$ws.Cells.Item(25, 2).Value = 19.72679021536619
$ws.Cells.Item(25, 3).Value = 6.59523365364207
$ws.Cells.Item(25, 4).Value = 4.739990692283265
$ws.Cells.Item(25, 5).Value = 10.635225877631
$ws.Cells.Item(25, 6).Value = 50.31876473595959
$ws.Cells.Item(25, 8).Value = 7.344005520526261
$ws.Cells.Item(25, 9).Value = 36.16035954471053
$ws.Cells.Item(25, 10).Value = 9.941912884130531
$ws.Cells.Item(25, 11).Value = 17.48570213262299
$ws.Cells.Item(25, 12).Value = 12.04650844682674
$ws.Cells.Item(25, 14).Value = 24.3682272436387
